# Apply updated cryptocurrency price and volume(1h) figures to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column cells to plain text so values like "26.813.38" or
# "301.14" are preserved exactly as strings, matching the source data feed.
$priceCells = @("D2", "D3", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D43", "D44", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.813.38"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "1.872.48"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "301.14"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "0.5320"
$ws.Range("E7").Value = "  +1.19%  "
$ws.Range("D8").Value = "0.3745"
$ws.Range("E8").Value = "  -1.92%  "
$ws.Range("D9").Value = "0.07185"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("D10").Value = "21.62"
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "0.8888"
$ws.Range("E11").Value = "  -1.85%  "
$ws.Range("D12").Value = "0.08174"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("D13").Value = "1.884.88"
$ws.Range("E13").Value = "  +13.58%  "
$ws.Range("D14").Value = "92.95"
$ws.Range("E14").Value = "  -3.66%  "
$ws.Range("D15").Value = "5.301"
$ws.Range("E15").Value = "  -1.34%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.07%  "
$ws.Range("D17").Value = "14.84"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("D18").Value = "0.000008526"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "26.864.83"
$ws.Range("D21").Value = "4.986"
$ws.Range("E21").Value = "  -2.76%  "
$ws.Range("D22").Value = "10.63"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").Value = "6.394"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").Value = "2.291"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").Value = "146.62"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("D26").Value = "1.742"
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "18.04"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("D28").Value = "114.01"
$ws.Range("E28").Value = "  -2.45%  "
$ws.Range("D29").Value = "4.712"
$ws.Range("E29").Value = "  -2.83%  "
$ws.Range("D30").Value = "4.611"
$ws.Range("E30").Value = "  -5.49%  "
$ws.Range("D31").Value = "0.09105"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").Value = "0.8107"
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "0.05017"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("D34").Value = "1.171"
$ws.Range("E34").Value = "  -4.70%  "
$ws.Range("D35").Value = "2.966"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "0.6097"
$ws.Range("E36").Value = "  +5.76%  "
$ws.Range("D37").Value = "2.661"
$ws.Range("E37").Value = "  -3.45%  "
$ws.Range("D38").Value = "3.201"
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("D39").Value = "0.01956"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "6.549"
$ws.Range("E41").Value = "  -0.99%  "
$ws.Range("D42").Value = "8.849"
$ws.Range("E42").Value = "  -2.84%  "
$ws.Range("D43").Value = "0.5180"
$ws.Range("E43").Value = "  +5.26%  "
$ws.Range("D44").Value = "115.04"
$ws.Range("E44").Value = "  -1.91%  "
$ws.Range("E45").Value = "  -1.86%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "1.640"
$ws.Range("E47").Value = "  -0.29%  "
$ws.Range("D48").Value = "9.947"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").Value = "37.51"
$ws.Range("E49").Value = "  -3.18%  "
$ws.Range("D50").Value = "0.06046"
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").Value = "62.34"
$ws.Range("E51").Value = "  -2.70%  "
